# ReportOrsha.xlsx update:
#  - rename sheet "Лист2" -> "Out"
#  - restore alphabetical sort (handled implicitly - sheet order unchanged)
#  - populate "Out" sheet with the "out-of-district manufacturer" certificate register table

$wb = $excel.ActiveWorkbook

# ---- 1. Rename the second sheet ----
$wsOut = $wb.Worksheets.Item(2)
$wsOut.Name = "Out"

# ---- 2. Column widths ----
$wsOut.Columns.Item(2).ColumnWidth = 31.875
$wsOut.Columns.Item(3).ColumnWidth = 11.09375
$wsOut.Columns.Item(4).ColumnWidth = 12.65625
$wsOut.Columns.Item(6).ColumnWidth = 12.65625
$wsOut.Columns.Item(7).ColumnWidth = 12.1875
$wsOut.Columns.Item(8).ColumnWidth = 18.125
$wsOut.Columns.Item(9).ColumnWidth = 45.46875

# ---- 3. Row heights ----
$wsOut.Rows.Item(1).RowHeight = 75.75
$wsOut.Rows.Item(3).RowHeight = 78.75

# ---- 4. Title row (A1:I1) ----
$title = "Реестр сертификатов продукции  собственного производства, выданных юридическим лицам и индивидуальным предпринимателям, зарегистрированным в Республике Беларусь с местом нахождения (жительства) в Оршанском районе и местом производства вне Оршанского района "
$titleRange = $wsOut.Range("A1:I1")
$titleRange.Merge()
$wsOut.Range("A1").Value = $title
$titleRange.Font.Name = "Times New Roman"
$titleRange.Font.Size = 14
$titleRange.HorizontalAlignment = -4108   # xlCenter
$titleRange.VerticalAlignment = -4160     # xlTop
$titleRange.WrapText = $true
$titleRange.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$titleRange.Borders.Item(9).Weight = 2      # xlThin

# ---- 5. Header row 2 + sub-header row 3 ----
$wsOut.Range("A2").Value = "№ п.п."
$wsOut.Range("B2").Value = "Производитель продукции                                   (место производства - вне Оршанского района)                         "
$wsOut.Range("C2").Value = "Дата выдачи сертификата"
$wsOut.Range("D2").Value = "Номера сертификата"
$wsOut.Range("E2").Value = "Срок действия сертификата"
$wsOut.Range("G2").Value = "Форма сертификата"
$wsOut.Range("H2").Value = "Укрупненное наименование продукции собственного производства"
$wsOut.Range("I2").Value = "Место нахождения производства: для продукции - на основании информации, указанной в сертификате; для работ, услуг - на основании информации в заявлении заказчика на экспертизу "

$wsOut.Range("E3").Value = "с"
$wsOut.Range("F3").Value = "до  "

# ---- 6. Merges for headers ----
$wsOut.Range("A2:A3").Merge()
$wsOut.Range("B2:B3").Merge()
$wsOut.Range("C2:C3").Merge()
$wsOut.Range("D2:D3").Merge()
$wsOut.Range("E2:F2").Merge()
$wsOut.Range("G2:G3").Merge()
$wsOut.Range("H2:H3").Merge()
$wsOut.Range("I2:I3").Merge()

# ---- 7. Formatting for header block (rows 2-3, columns A-I) ----
$headerBlock = $wsOut.Range("A2:I3")
$headerBlock.Font.Name = "Times New Roman"
$headerBlock.Font.Size = 11
$headerBlock.VerticalAlignment = -4160   # xlTop
$headerBlock.HorizontalAlignment = -4108 # xlCenter

# Wrap text for the cells that wrap (all except A and the "с"/"до" single-word cells)
$wsOut.Range("B2:B3").HorizontalAlignment = -4131  # xlLeft
$wsOut.Range("I2:I3").HorizontalAlignment = -4131  # xlLeft
$wsOut.Range("B2:B3,C2:C3,D2:D3,E2:F2,E3,F3,G2:G3,H2:H3,I2:I3").WrapText = $true

# A2:A3 and E3/F3 do not wrap (single short values)
$wsOut.Range("A2:A3").WrapText = $false

# ---- 8. Borders: thin box around the whole header block, plus internal gridlines ----
$outline = $wsOut.Range("A2:I3")
$outline.Borders.Item(7).LineStyle = 1  # xlEdgeLeft
$outline.Borders.Item(7).Weight = 2
$outline.Borders.Item(10).LineStyle = 1 # xlEdgeRight
$outline.Borders.Item(10).Weight = 2
$outline.Borders.Item(8).LineStyle = 1  # xlEdgeTop
$outline.Borders.Item(8).Weight = 2
$outline.Borders.Item(9).LineStyle = 1  # xlEdgeBottom
$outline.Borders.Item(9).Weight = 2
$outline.Borders.Item(11).LineStyle = 1 # xlInsideVertical
$outline.Borders.Item(11).Weight = 2
$outline.Borders.Item(12).LineStyle = 1 # xlInsideHorizontal
$outline.Borders.Item(12).Weight = 2

# ---- 9. Page setup ----
$wsOut.PageSetup.PaperSize = 9     # xlPaperA4
$wsOut.PageSetup.Orientation = 1   # xlPortrait

Write-Host "Out sheet populated"
